$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.459.51"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "2.163.20"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  -0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "228.91"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +1.25%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "63.29"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +5.08%  "
$ws.Range("E8").Value = "  -0.06%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.397"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.23%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0867"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +7.19%  "
$ws.Range("D13").Value = "2.483.20"
$ws.Range("E13").Value = "  +3.44%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "22.23"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.58%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.818"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D17").Value = "2.172.61"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "39.397.98"
$ws.Range("E18").Value = "  +1.83%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "72.38"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +1.85%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "228.70"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E24").Value = "  -1.49%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.36"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.78"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.76%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "171.93"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  +8.02%  "
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +4.06%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.81"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.40%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "7.07"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +9.18%  "
$ws.Range("E36").Value = "  +2.32%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.43"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -0.11%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "18.18"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "103.94"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "1.532.41"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  +5.61%  "
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("E46").Value = "  +7.21%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.79"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "4.24"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("D50").Value = "2.366.81"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("E51").Value = "  +0.37%  "
